$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Candidate ID numbers in column B
$ws.Range("B2").Value = 231102244
$ws.Range("B3").Value = 231102243
$ws.Range("B4").Value = 231102242
$ws.Range("B5").Value = 231102241
$ws.Range("B6").Value = 231102240

# Row 2 (Candidate ID 231102244)
$ws.Range("A2").Value = "JuwWN785"
$ws.Range("C2").Value = "tbxreyl54"
$ws.Range("D2").Value = "UF6#kj9%"
$ws.Range("F2").Value = "iZmSWKXL"
$ws.Range("G2").Value = "XXRA"

# Row 3 (Candidate ID 231102243)
$ws.Range("A3").Value = "ilAbq356"
$ws.Range("C3").Value = "psulted40"
$ws.Range("D3").Value = "Ta`$bU7!3"
$ws.Range("F3").Value = "BBrsqHli"
$ws.Range("G3").Value = "dmHP"

# Row 4 (Candidate ID 231102242)
$ws.Range("A4").Value = "QxsZx405"
$ws.Range("C4").Value = "mftatws42"
$ws.Range("D4").Value = "x&29RyZ!"
$ws.Range("F4").Value = "TqDKIpSs"
$ws.Range("G4").Value = "wTns"

# Row 5 (Candidate ID 231102241)
$ws.Range("A5").Value = "yNbqA358"
$ws.Range("C5").Value = "ynzlilw95"
$ws.Range("D5").Value = "Cc#6B&u4"
$ws.Range("F5").Value = "iLqXOXwy"
$ws.Range("G5").Value = "FHQf"

# Row 6 (Candidate ID 231102240)
$ws.Range("A6").Value = "OeZwP266"
$ws.Range("C6").Value = "dnforsk61"
$ws.Range("D6").Value = "p!tZ4#9G"
$ws.Range("F6").Value = "bDgDhhOt"
$ws.Range("G6").Value = "qsrc"
